# edit.ps1
# Applies the "Updated cryptos list" diff: refresh Price (D) and Volume(1h) (E)
# values for each coin row, plus the swap of the NEARProtocol/Bittensor rows
# (rows 33 and 34 exchange places with updated figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain TEXT so Excel does not silently reinterpret
# digit/period strings (e.g. "598.49") as numbers. We temporarily force a
# text number-format, assign the value, then restore the cell's original
# ("Normal") style so no stray formatting is left on the cell.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '70.555.73'
$ws.Range('E2').Value = '  +2.32%  '
$ws.Range('D3').Value = '3.563.81'
$ws.Range('E3').Value = '  +1.21%  '
$ws.Range('E4').Value = '  +0.14%  '
Set-TextValue $ws.Range('D5') '598.49'
$ws.Range('E5').Value = '  +1.71%  '
Set-TextValue $ws.Range('D6') '172.33'
$ws.Range('E6').Value = '  +1.37%  '
$ws.Range('D7').Value = '3.557.80'
$ws.Range('E7').Value = '  +1.19%  '
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('E9').Value = '  +0.02%  '
Set-TextValue $ws.Range('D10') '0.197'
$ws.Range('E10').Value = '  +5.00%  '
Set-TextValue $ws.Range('D11') '7.40'
$ws.Range('E11').Value = '  +8.69%  '
Set-TextValue $ws.Range('D12') '0.587'
$ws.Range('E12').Value = '  +1.49%  '
Set-TextValue $ws.Range('D13') '46.35'
$ws.Range('E13').Value = '  -2.01%  '
$ws.Range('E14').Value = '  +0.47%  '
$ws.Range('D15').Value = '4.137.96'
$ws.Range('E15').Value = '  +1.24%  '
Set-TextValue $ws.Range('D16') '8.37'
$ws.Range('E16').Value = '  -0.82%  '
Set-TextValue $ws.Range('D17') '610.40'
$ws.Range('E17').Value = '  -0.95%  '
$ws.Range('D18').Value = '3.574.47'
$ws.Range('E18').Value = '  +1.89%  '
$ws.Range('D19').Value = '70.600.98'
$ws.Range('E19').Value = '  +2.43%  '
$ws.Range('E20').Value = '  -1.05%  '
Set-TextValue $ws.Range('D21') '17.37'
$ws.Range('E21').Value = '  -0.56%  '
Set-TextValue $ws.Range('D22') '0.882'
$ws.Range('E22').Value = '  -0.36%  '
Set-TextValue $ws.Range('D23') '9.23'
$ws.Range('E23').Value = '  -17.19%  '
Set-TextValue $ws.Range('D24') '15.78'
$ws.Range('E24').Value = '  +0.14%  '
Set-TextValue $ws.Range('D25') '96.70'
$ws.Range('E25').Value = '  +0.14%  '
Set-TextValue $ws.Range('D26') '3.75'
$ws.Range('E26').Value = '  -2.03%  '
Set-TextValue $ws.Range('D27') '1.00'
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('E28').Value = '  +0.12%  '
Set-TextValue $ws.Range('D29') '33.78'
$ws.Range('E29').Value = '  +3.63%  '
Set-TextValue $ws.Range('D30') '9.09'
$ws.Range('E30').Value = '  -1.55%  '
Set-TextValue $ws.Range('D31') '8.31'
$ws.Range('E31').Value = '  -2.42%  '
$ws.Range('E32').Value = '  -2.47%  '
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D33') '7.12'
$ws.Range('E33').Value = '  +2.86%  '
$ws.Range('B34').Value = 'Bittensor'
$ws.Range('C34').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range('D34') '653.27'
$ws.Range('E34').Value = '  +4.75%  '
Set-TextValue $ws.Range('D35') '1.29'
$ws.Range('E35').Value = '  -1.25%  '
Set-TextValue $ws.Range('D36') '3.64'
$ws.Range('E36').Value = '  +6.14%  '
$ws.Range('E37').Value = '  -1.43%  '
Set-TextValue $ws.Range('D38') '10.78'
$ws.Range('E38').Value = '  +0.36%  '
Set-TextValue $ws.Range('D39') '0.0474'
$ws.Range('E39').Value = '  +6.73%  '
Set-TextValue $ws.Range('D40') '57.18'
$ws.Range('E40').Value = '  +0.09%  '
$ws.Range('E41').Value = '  +0.27%  '
$ws.Range('E42').Value = '  +5.15%  '
$ws.Range('D43').Value = '3.382.23'
$ws.Range('E43').Value = '  -0.03%  '
Set-TextValue $ws.Range('D44') '0.320'
$ws.Range('E44').Value = '  -1.79%  '
$ws.Range('D45').Value = '0.0₃0708'
$ws.Range('E45').Value = '  +1.98%  '
Set-TextValue $ws.Range('D46') '32.78'
$ws.Range('E46').Value = '  +0.01%  '
Set-TextValue $ws.Range('D47') '2.95'
$ws.Range('E47').Value = '  +6.94%  '
$ws.Range('E48').Value = '  +4.80%  '
$ws.Range('E49').Value = '  +0.66%  '
Set-TextValue $ws.Range('D50') '132.37'
$ws.Range('E50').Value = '  -0.63%  '
$ws.Range('E51').Value = '  -0.05%  '
